$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new columns before column C (old C shifts to E)
$ws.Range("C1:D1").EntireColumn.Insert()

# The old B1 header ("Jun_13") needs to move to D1 before we overwrite B1
$ws.Range("D1").Value = $ws.Range("B1").Value2

# New header values
$ws.Range("B1").Value = "Jun_17"
$ws.Range("C1").Value = "Jun_15"

# Fill the new C,D columns (rows 2-27) with "UN" to match column B
$ws.Range("C2:D27").Value = "UN"

# Match column width of the neighboring date columns
$ws.Range("C:D").ColumnWidth = 57.75
